$d = $word.ActiveDocument

# The document currently has a single empty paragraph. We need to end up with:
#   p1: "2.1"
#   p2: big multi-run paragraph about the research study
#   p3: (the original, untouched, empty paragraph)
#   p4: "2.2"

# Step 1: push the original empty paragraph down by inserting a placeholder
# paragraph for "2.1" right before it.
$target = $d.Paragraphs(1).Range
$target.InsertBefore("2.1`r")

# Step 2: the original empty paragraph is now Paragraphs(2). Insert another
# placeholder paragraph before it, to be filled in with the big paragraph's
# multiple runs.
$target2 = $d.Paragraphs(2).Range
$target2.InsertBefore("PLACEHOLDER`r")

# Step 3: fill the placeholder paragraph (now Paragraphs(2)) with the real
# multi-run content via InsertXML so the run boundaries match exactly.
$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$body = '<w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">The research paper "Goal Statistics and Scoring Attributes of the 2018 FIFA World Cup" by Ankur Biswas and Nita Bandyopadhyay analyses the patterns and trends of goals scored during the 2018 FIFA World Cup in Russia. This research was conducted using data obtained from the official FIFA database. It focused on key variables such as goals, shots, shots on target, ball possession, passes, pass accuracy, fouls, yellow cards, offsides and corners. The study found that the average number of goals per match was 2.64, with the majority scored by the winning teams. There were significant differences between the winning and losing teams in terms of shots, shots on target and yellow cards. </w:t></w:r>' + `
  '<w:r><w:t>Also,</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> this study highlighted significant correlations between most of the variables, providing valuable insights for coaches and players to develop effective game plans and training schedules</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> ()</w:t></w:r>' + `
  '<w:r><w:t>.</w:t></w:r>' + `
  '</w:p></w:body>'

$p2 = $d.Paragraphs(2).Range
[void]$p2.InsertXML($xmlHeader + $body + $xmlFooter)

# Step 4: append the "2.2" paragraph right after the (still untouched) empty
# paragraph, which is now Paragraphs(3).
$emptyPara = $d.Paragraphs(3).Range
$emptyPara.InsertAfter("`r2.2")
